# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45903
$ws.Range("B2").Value = 105.23
$ws.Range("D2").Value = 85
$ws.Range("E2").Value = 80.5
$ws.Range("F2").Value = 76.05
$ws.Range("G2").Value = 76.05
$ws.Range("H2").Value = 83.22
$ws.Range("I2").Value = 93
$ws.Range("J2").Value = 83.22
$ws.Range("K2").Value = 17.5
$ws.Range("L2").Value = 2.01
$ws.Range("M2").Value = 2.01
$ws.Range("N2").Value = 0.5
$ws.Range("O2").Value = 0.5
$ws.Range("P2").Value = 0.05
$ws.Range("Q2").Value = 0.05
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0.05
$ws.Range("T2").Value = 12.85
$ws.Range("U2").Value = 70
$ws.Range("V2").Value = 103.13
$ws.Range("W2").Value = 133.45
$ws.Range("X2").Value = 107.43
$ws.Range("Y2").Value = 100.86
$ws.Range("Z2").Value = 55.32
$ws.Range("AB2").Value = 111.22
$ws.Range("AD2").Value = 118.29
$ws.Range("AF2").Value = 104.15
